$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $ws.Range("ZZ1").NumberFormat = "@"
    $ws.Range("ZZ1").Value = $val
    $ws.Range("ZZ1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4104)
    $ws.Range("ZZ1").Clear()
}

Set-TextValue 'D2' '29.405.57'
Set-TextValue 'E2' '  +0.09%  '
Set-TextValue 'D3' '1.843.46'
Set-TextValue 'E3' '  -0.03%  '
Set-TextValue 'D4' '0.9984'
Set-TextValue 'E4' '  -0.32%  '
Set-TextValue 'D5' '240.23'
Set-TextValue 'E5' '  -0.13%  '
Set-TextValue 'D6' '0.6325'
Set-TextValue 'E6' '  +0.83%  '
Set-TextValue 'E7' '  -0.23%  '
Set-TextValue 'D8' '0.07474'
Set-TextValue 'E8' '  -0.11%  '
Set-TextValue 'D9' '25.10'
Set-TextValue 'E9' '  +3.04%  '
Set-TextValue 'D10' '0.2907'
Set-TextValue 'E10' '  +0.45%  '
Set-TextValue 'D11' '0.07742'
Set-TextValue 'E11' '  +0.21%  '
Set-TextValue 'D12' '1.848.18'
Set-TextValue 'E12' '  +0.25%  '
Set-TextValue 'D13' '4.991'
Set-TextValue 'E13' '  +0.14%  '
Set-TextValue 'D14' '0.6792'
Set-TextValue 'E14' '  +0.10%  '
Set-TextValue 'E15' '  -0.10%  '
Set-TextValue 'D16' '82.14'
Set-TextValue 'E16' '  +0.02%  '
Set-TextValue 'D17' '6.277'
Set-TextValue 'E17' '  +3.05%  '
Set-TextValue 'D18' '29.474.41'
Set-TextValue 'E18' '  +0.21%  '
Set-TextValue 'D19' '229.87'
Set-TextValue 'E19' '  +0.54%  '
Set-TextValue 'E20' '  +0.62%  '
Set-TextValue 'D21' '0.9994'
Set-TextValue 'E21' '  -0.28%  '
Set-TextValue 'D22' '7.440'
Set-TextValue 'E22' '  +1.01%  '
Set-TextValue 'E23' '  -0.14%  '
Set-TextValue 'D24' '158.43'
Set-TextValue 'E24' '  -0.30%  '
Set-TextValue 'D25' '8.505'
Set-TextValue 'E25' '  +1.39%  '
Set-TextValue 'E26' '  -1.66%  '
Set-TextValue 'E27' '  -0.38%  '
Set-TextValue 'D28' '0.06595'
Set-TextValue 'E28' '  +16.06%  '
Set-TextValue 'D29' '1.428'
Set-TextValue 'E29' '  +2.55%  '
Set-TextValue 'E30' '  +0.56%  '
Set-TextValue 'D31' '4.079'
Set-TextValue 'E31' '  -0.51%  '
Set-TextValue 'D32' '4.058'
Set-TextValue 'E32' '  +0.44%  '
Set-TextValue 'D33' '1.841'
Set-TextValue 'E34' '  -0.30%  '
Set-TextValue 'D35' '0.6991'
Set-TextValue 'E35' '  +1.23%  '
Set-TextValue 'D37' '0.01860'
Set-TextValue 'E37' '  +2.45%  '
Set-TextValue 'D38' '1.249.32'
Set-TextValue 'E38' '  -0.18%  '
Set-TextValue 'D39' '2.816'
Set-TextValue 'E39' '  -1.33%  '
Set-TextValue 'D40' '6.795'
Set-TextValue 'E40' '  +4.27%  '
Set-TextValue 'D41' '0.9344'
Set-TextValue 'E41' '  +3.32%  '
Set-TextValue 'D42' '0.9994'
Set-TextValue 'E42' '  -0.16%  '
Set-TextValue 'D43' '2.018.77'
Set-TextValue 'E43' '  +0.71%  '
Set-TextValue 'D44' '101.07'
Set-TextValue 'E44' '  -0.12%  '
Set-TextValue 'D45' '65.49'
Set-TextValue 'E45' '  -0.23%  '
Set-TextValue 'B46' 'BabyDogeCoin'
Set-TextValue 'C46' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D46' '0.00000000119'
Set-TextValue 'E46' '  +3.83%  '
Set-TextValue 'B47' 'Aptos'
Set-TextValue 'C47' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D47' '7.070'
Set-TextValue 'E47' '  -0.09%  '
Set-TextValue 'B48' 'RenderToken'
Set-TextValue 'C48' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D48' '1.718'
Set-TextValue 'E48' '  +3.87%  '
Set-TextValue 'D49' '9.059'
Set-TextValue 'E49' '  +1.03%  '
Set-TextValue 'D50' '0.1150'
Set-TextValue 'E50' '  -0.89%  '
Set-TextValue 'D51' '0.3912'
Set-TextValue 'E51' '  -0.65%  '
